$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Customer + Account rows (TC, Customer_ID, PD) to append below the
# existing table (which currently ends at row 25).
$newRows = @(
    @("118500", "17704735", "6004"),
    @("118518", "17704736", "6020"),
    @("118498", "17704737", "1001"),
    @("118452", "17704738", "1001"),
    @("118518", "17704739", "6012"),
    @("118448", "17704740", "1047"),
    @("118518", "17704741", "1035"),
    @("118452", "17704742", "1150"),
    @("118448", "17704743", "1068"),
    @("118448", "17704745", "1005"),
    @("118500", "17704746", "6004"),
    @("118500", "17704747", "6004"),
    @("118500", "17704748", "6004"),
    @("118518", "17704749", "6020"),
    @("118498", "17704750", "1001"),
    @("118452", "17704751", "1001"),
    @("118518", "17704752", "6012"),
    @("118448", "17704753", "1047"),
    @("118518", "17704754", "1035"),
    @("118452", "17704755", "1150"),
    @("118448", "17704756", "1068"),
    @("118448", "17704757", "1005")
)

$startRow = 26
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $vals = $newRows[$i]
    # Leading "'" forces each value to be stored as text (matching the
    # source data, which is all numeric-looking strings such as customer
    # and account ids) instead of being auto-coerced to a number.
    $ws.Range("A$r").Formula = "'" + $vals[0]
    $ws.Range("B$r").Formula = "'" + $vals[1]
    $ws.Range("C$r").Formula = "'" + $vals[2]
}

# Drop the quote-prefix formatting flag picked up above so the new cells
# keep the sheet's default (unstyled) appearance.
$ws.Range("A26:C47").ClearFormats()
